# Daily attendance processing - 2025-12-31 06:44:01
# Swap the order of names in the "Recorded By" column (G) from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
